# Generate Report for Handback
$wb = $excel.ActiveWorkbook

# --- Overview sheet: status text for the 61aaaf58 md file changes ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Cells.Item(3, 5).Value = "Handback transform failed"
$overview.Cells.Item(3, 6).Value = "Handback transform failed"

# --- zh-cn sheet: widen "Error Detail" column (P), update Status (C) and set error detail for row 3 ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Columns.Item(16).ColumnWidth = 39.166666666666664
$zhcn.Cells.Item(3, 3).Value = "Handback transform failed"
$zhcn.Cells.Item(3, 16).Value = "Handback file name: nzhntt4u.3fy is different with handoff file name: 61aaaf58-1911-40d7-9bf0-12930459d975.d77335cf432a548b305037f7565b4d22b62c380a.zh-cn."

# --- de-de sheet: widen "Error Detail" column (P), update Status (C) and set error detail for row 3 ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Columns.Item(16).ColumnWidth = 39.166666666666664
$dede.Cells.Item(3, 3).Value = "Handback transform failed"
$dede.Cells.Item(3, 16).Value = "Handback file name: nzhntt4u.3fy is different with handoff file name: 61aaaf58-1911-40d7-9bf0-12930459d975.d77335cf432a548b305037f7565b4d22b62c380a.de-de."
